$d = $word.ActiveDocument

# --- Locate the sentence to edit -----------------------------------------
# Original run text:
#   "Decision percentages: how many food pod, run pod entries and how many
#    exits from T-maze?"
# Target (split across three runs):
#   "Decision percentages" + " and raw numbers" + ": how many food pod, ..."
#
# We find the phrase "Decision percentages" and collapse the range to the
# point right after it (i.e. right before the colon), then insert the new
# text " and raw numbers" there.

$findRange = $d.Content
$found = $findRange.Find.Execute("Decision percentages", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Decision percentages' in the document."
}

# Collapse to the end of the match (right before ": how many food pod ...").
$findRange.Collapse(0)

# Use tracked-changes insertion so the new text lands in its own Run
# (matching the diff's 3-run split) instead of being silently merged back
# into the neighbouring run. We then accept the single resulting revision,
# which leaves plain, un-marked-up runs behind.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$findRange.InsertAfter(" and raw numbers")

$d.TrackRevisions = $wasTracking

while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
